# Weekly update: insert a new "Albahaca" price observation as the new
# first data row (row 102) for "Terminal La Palmera de La Serena", and
# push every following observation down by one row (the previous last
# row, 183, becomes row 184).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 102:183 down to 103:184 by inserting a blank row at 102.
$ws.Rows("102:102").Insert()

# Populate the newly inserted row 102 with the new weekly observation.
$ws.Range("A102").Value = 8
$ws.Range("B102").Value = "Terminal La Palmera de La Serena"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 45072
$ws.Range("E102").Value = 4
$ws.Range("F102").Value = 100112052
$ws.Range("G102").Value = "Albahaca"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 2800
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 2900
$ws.Range("N102").Value = "$/paquete"
$ws.Range("O102").Value = "Región de Arica y Parinacota"
$ws.Range("P102").Value = 2900
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = "Hortaliza"
